# Nexial "#system" lookup-table maintenance:
# add outputToCloud(resource) to the `base` command list, add a new
# `text` command group (spellCheck) and register it in the `target`
# drop-down list, shifting the existing alphabetical entries down/right
# to keep everything sorted.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ------------------------------------------------------------------
# 1) "base" column (E): insert "outputToCloud(resource)" before
#    "prependText(var,prependWith)" -> new E22, shifting E22:E38 down
#    to E23:E39.
# ------------------------------------------------------------------
for ($r = 38; $r -ge 22; $r--) {
    $v = $ws.Cells.Item($r, 5).Value()
    $ws.Cells.Item($r + 1, 5).Value = $v
}
$ws.Cells.Item(22, 5).Value = "outputToCloud(resource)"

# ------------------------------------------------------------------
# 2) "target" column (A): insert "text" before "web" -> new A25,
#    shifting A25:A30 down to A26:A31.
# ------------------------------------------------------------------
for ($r = 30; $r -ge 25; $r--) {
    $v = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r + 1, 1).Value = $v
}
$ws.Cells.Item(25, 1).Value = "text"

# ------------------------------------------------------------------
# 3) Insert a brand-new column before Y (col 25) for the new `text`
#    command group. Everything from Y onward (web, webalert,
#    webcookie, ws, ws.async, xml) shifts one column to the right.
# ------------------------------------------------------------------
$ws.Columns.Item(25).Insert()
$ws.Cells.Item(1, 25).Value = "text"
$ws.Cells.Item(2, 25).Value = "spellCheck(var,profile,text)"

# ------------------------------------------------------------------
# 4) Update the named ranges to reflect the new extents/columns.
# ------------------------------------------------------------------
$wb.Names.Item("base").RefersTo = "='#system'!`$E`$2:`$E`$39"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$31"
$wb.Names.Item("web").RefersTo = "='#system'!`$Z`$2:`$Z`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AC`$2:`$AC`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AD`$2:`$AD`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AE`$2:`$AE`$27"
$wb.Names.Add("text", "='#system'!`$Y`$2:`$Y`$2")
